$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reflect the final selection state (user had selected A2:E2 before editing)
[void]$ws.Range("A2:E2").Select()

# B2: password for "kikaho" account was changed
$ws.Range("B2").Value = "Nhacc123@"

# C2, D2, E2: stray accounts ("abc"/"123"/"zcx") were wiped out, but the
# cells themselves were left behind holding empty text (quote-prefix trick
# forces Excel to keep a real, empty shared-string cell instead of
# deleting it outright), then restore default "Normal" styling.
$ws.Range("C2").Value = "'"
$ws.Range("C2").Style = "Normal"

$ws.Range("D2").Value = "'"
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").Value = "'"
$ws.Range("E2").Style = "Normal"
